$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text, matching the
# original inline-string cell contents (e.g. "0.9998", "313.17").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.328.36"
$ws.Range("E2").Value = "  -0.64%  "

$ws.Range("D3").Value = "1.812.38"
$ws.Range("E3").Value = "  -0.80%  "

$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.50%  "

$ws.Range("D5").Value = "313.17"
$ws.Range("E5").Value = "  -1.09%  "

$ws.Range("D6").Value = "0.9992"
$ws.Range("E6").Value = "  -0.50%  "

$ws.Range("D7").Value = "0.5161"

$ws.Range("D8").Value = "0.4035"
$ws.Range("E8").Value = "  +4.42%  "

$ws.Range("D9").Value = "0.07875"
$ws.Range("E9").Value = "  -5.39%  "

$ws.Range("D10").Value = "1.115"
$ws.Range("E10").Value = "  -0.47%  "

$ws.Range("D11").Value = "41.12"
$ws.Range("E11").Value = "  -2.20%  "

$ws.Range("D12").Value = "6.372"
$ws.Range("E12").Value = "  -0.60%  "

$ws.Range("D13").Value = "0.9994"
$ws.Range("E13").Value = "  -0.50%  "

$ws.Range("D14").Value = "20.49"
$ws.Range("E14").Value = "  -3.17%  "

$ws.Range("D15").Value = "7.355"
$ws.Range("E15").Value = "  -1.95%  "

$ws.Range("D16").Value = "1.799.16"
$ws.Range("E16").Value = "  -1.46%  "

$ws.Range("D17").Value = "92.91"

$ws.Range("E18").Value = "  -3.66%  "

$ws.Range("D19").Value = "0.06595"
$ws.Range("E19").Value = "  -0.79%  "

$ws.Range("D20").Value = "0.9986"
$ws.Range("E20").Value = "  -0.56%  "

$ws.Range("D21").Value = "17.37"
$ws.Range("E21").Value = "  -2.34%  "

$ws.Range("D22").Value = "6.058"
$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("D23").Value = "28.381.89"
$ws.Range("E23").Value = "  -0.60%  "

$ws.Range("D24").Value = "11.21"
$ws.Range("E24").Value = "  -1.77%  "

$ws.Range("D25").Value = "2.226"
$ws.Range("E25").Value = "  -2.78%  "

$ws.Range("D26").Value = "160.78"
$ws.Range("E26").Value = "  +0.66%  "

$ws.Range("D27").Value = "20.65"
$ws.Range("E27").Value = "  -2.56%  "

$ws.Range("D28").Value = "2.015.04"
$ws.Range("E28").Value = "  -0.94%  "

$ws.Range("D29").Value = "2.412"
$ws.Range("E29").Value = "  +0.31%  "

$ws.Range("D30").Value = "128.84"
$ws.Range("E30").Value = "  +2.22%  "

$ws.Range("D31").Value = "0.1086"
$ws.Range("E31").Value = "  -0.60%  "

$ws.Range("D32").Value = "1.055"
$ws.Range("E32").Value = "  -3.71%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "5.597"
$ws.Range("E33").Value = "  -2.37%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "3.659"
$ws.Range("E34").Value = "  -0.55%  "

$ws.Range("D35").Value = "0.07199"
$ws.Range("E35").Value = "  -5.44%  "

$ws.Range("D36").Value = "9.159"
$ws.Range("E36").Value = "  +4.67%  "

$ws.Range("D37").Value = "0.02335"
$ws.Range("E37").Value = "  -1.63%  "

$ws.Range("D38").Value = "0.2168"
$ws.Range("E38").Value = "  -2.68%  "

$ws.Range("E39").Value = "  -0.62%  "

$ws.Range("D40").Value = "5.075"
$ws.Range("E40").Value = "  -3.61%  "

$ws.Range("D41").Value = "0.6241"
$ws.Range("E41").Value = "  -2.14%  "

$ws.Range("D42").Value = "0.9985"
$ws.Range("E42").Value = "  -0.56%  "

$ws.Range("D43").Value = "1.158"
$ws.Range("E43").Value = "  -2.68%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "13.30"
$ws.Range("E44").Value = "  -2.13%  "

$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "1.323"
$ws.Range("E45").Value = "  -5.35%  "

$ws.Range("D46").Value = "0.6023"
$ws.Range("E46").Value = "  -1.36%  "

$ws.Range("D47").Value = "3.746"
$ws.Range("E47").Value = "  -1.38%  "

$ws.Range("D48").Value = "126.28"
$ws.Range("E48").Value = "  -0.94%  "

$ws.Range("D49").Value = "1.218"
$ws.Range("E49").Value = "  +0.87%  "

$ws.Range("D50").Value = "1.944"

$ws.Range("E51").Value = "  -1.78%  "
